$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item(1)

# Shorten responsible-person names in the RESPONSÁVEL column (C)
$ws.Range("C6").Value = "Pedro H e Pedro A"
$ws.Range("C3").Value = "Pedro H, Pedro A, Douglas"
$ws.Range("C4").Value = "Pedro H, Pedro A, Douglas"
$ws.Range("C5").Value = "Pedro H, Pedro A, Douglas"

# Leave the last active selection on F7, matching the saved file state
$ws.Range("F7").Select()
